$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner in A1
$ws.Range("A1").Value = 'Datos actualizados a 23 de Octubre de 2020 a las 17:13'

# Rebuild the country data table (rows 4-221), sorted by Casos totales descending,
# reflecting refreshed case counts and re-ranked countries (e.g. Portugal, Reunion,
# San Marino, Montserrat moving up in rank).
$nrows = 218
$ncols = 8
$data = New-Object 'object[,]' $nrows,$ncols

$data[0,0] = 'Estados Unidos'
$data[0,1] = 8675170
$data[0,2] = 13519
$data[0,3] = 5662993
$data[0,4] = 2783600
$data[0,5] = 0
$data[0,6] = 196
$data[0,7] = 228577
$data[1,0] = 'India'
$data[1,1] = 7772649
$data[1,2] = 13009
$data[1,3] = 6958897
$data[1,4] = 696316
$data[1,5] = 0
$data[1,6] = 100
$data[1,7] = 117436
$data[2,0] = 'Brasil'
$data[2,1] = 5332634
$data[2,2] = 0
$data[2,3] = 4785297
$data[2,4] = 391375
$data[2,5] = 0
$data[2,6] = 0
$data[2,7] = 155962
$data[3,0] = 'Rusia'
$data[3,1] = 1480646
$data[3,2] = 17340
$data[3,3] = 1119251
$data[3,4] = 335870
$data[3,5] = 0
$data[3,6] = 283
$data[3,7] = 25525
$data[4,0] = 'España'
$data[4,1] = 1090521
$data[4,2] = 0
$data[4,3] = 0
$data[4,4] = 0
$data[4,5] = 0
$data[4,6] = 0
$data[4,7] = 34521
$data[5,0] = 'Argentina'
$data[5,1] = 1053650
$data[5,2] = 0
$data[5,3] = 851854
$data[5,4] = 173839
$data[5,5] = 0
$data[5,6] = 0
$data[5,7] = 27957
$data[6,0] = 'Francia'
$data[6,1] = 999043
$data[6,2] = 0
$data[6,3] = 108599
$data[6,4] = 856234
$data[6,5] = 0
$data[6,6] = 0
$data[6,7] = 34210
$data[7,0] = 'Colombia'
$data[7,1] = 990270
$data[7,2] = 0
$data[7,3] = 893712
$data[7,4] = 66922
$data[7,5] = 0
$data[7,6] = 0
$data[7,7] = 29636
$data[8,0] = 'Peru'
$data[8,1] = 879876
$data[8,2] = 0
$data[8,3] = 796719
$data[8,4] = 49173
$data[8,5] = 0
$data[8,6] = 0
$data[8,7] = 33984
$data[9,0] = 'Mexico'
$data[9,1] = 874171
$data[9,2] = 6612
$data[9,3] = 636391
$data[9,4] = 149886
$data[9,5] = 0
$data[9,6] = 479
$data[9,7] = 87894
$data[10,0] = 'Reino Unido'
$data[10,1] = 810467
$data[10,2] = 0
$data[10,3] = 0
$data[10,4] = 0
$data[10,5] = 0
$data[10,6] = 0
$data[10,7] = 44347
$data[11,0] = 'Sudafrica'
$data[11,1] = 710515
$data[11,2] = 0
$data[11,3] = 642560
$data[11,4] = 49112
$data[11,5] = 0
$data[11,6] = 0
$data[11,7] = 18843
$data[12,0] = 'Iran'
$data[12,1] = 556891
$data[12,2] = 6134
$data[12,3] = 446685
$data[12,4] = 78221
$data[12,5] = 0
$data[12,6] = 335
$data[12,7] = 31985
$data[13,0] = 'Chile'
$data[13,1] = 498906
$data[13,2] = 1775
$data[13,3] = 471343
$data[13,4] = 13719
$data[13,5] = 0
$data[13,6] = 52
$data[13,7] = 13844
$data[14,0] = 'Italia'
$data[14,1] = 484869
$data[14,2] = 19143
$data[14,3] = 261808
$data[14,4] = 186002
$data[14,5] = 0
$data[14,6] = 91
$data[14,7] = 37059
$data[15,0] = 'Irak'
$data[15,1] = 445949
$data[15,2] = 3785
$data[15,3] = 375188
$data[15,4] = 60248
$data[15,5] = 0
$data[15,6] = 48
$data[15,7] = 10513
$data[16,0] = 'Alemania'
$data[16,1] = 413139
$data[16,2] = 9265
$data[16,3] = 310200
$data[16,4] = 92870
$data[16,5] = 0
$data[16,6] = 25
$data[16,7] = 10069
$data[17,0] = 'Banglades'
$data[17,1] = 396413
$data[17,2] = 1586
$data[17,3] = 312065
$data[17,4] = 78587
$data[17,5] = 0
$data[17,6] = 14
$data[17,7] = 5761
$data[18,0] = 'Indonesia'
$data[18,1] = 381910
$data[18,2] = 4369
$data[18,3] = 305100
$data[18,4] = 63733
$data[18,5] = 0
$data[18,6] = 118
$data[18,7] = 13077
$data[19,0] = 'Filipinas'
$data[19,1] = 365799
$data[19,2] = 1923
$data[19,3] = 312691
$data[19,4] = 46193
$data[19,5] = 0
$data[19,6] = 132
$data[19,7] = 6915
$data[20,0] = 'Turquia'
$data[20,1] = 355528
$data[20,2] = 0
$data[20,3] = 310027
$data[20,4] = 35917
$data[20,5] = 0
$data[20,6] = 0
$data[20,7] = 9584
$data[21,0] = 'Arabia Saudita'
$data[21,1] = 344157
$data[21,2] = 383
$data[21,3] = 330578
$data[21,4] = 8315
$data[21,5] = 0
$data[21,6] = 14
$data[21,7] = 5264
$data[22,0] = 'Ucrania'
$data[22,1] = 330396
$data[22,2] = 7517
$data[22,3] = 137578
$data[22,4] = 186654
$data[22,5] = 0
$data[22,6] = 121
$data[22,7] = 6164
$data[23,0] = 'Pakistan'
$data[23,1] = 326216
$data[23,2] = 736
$data[23,3] = 309646
$data[23,4] = 9855
$data[23,5] = 0
$data[23,6] = 13
$data[23,7] = 6715
$data[24,0] = 'Israel'
$data[24,1] = 308840
$data[24,2] = 593
$data[24,3] = 290130
$data[24,4] = 16381
$data[24,5] = 0
$data[24,6] = 10
$data[24,7] = 2329
$data[25,0] = 'Paises Bajos'
$data[25,1] = 272401
$data[25,2] = 9996
$data[25,3] = 0
$data[25,4] = 0
$data[25,5] = 0
$data[25,6] = 45
$data[25,7] = 6964
$data[26,0] = 'Belgica'
$data[26,1] = 270132
$data[26,2] = 16746
$data[26,3] = 22213
$data[26,4] = 237331
$data[26,5] = 0
$data[26,6] = 49
$data[26,7] = 10588
$data[27,0] = 'Polonia'
$data[27,1] = 228318
$data[27,2] = 13632
$data[27,3] = 105092
$data[27,4] = 119054
$data[27,5] = 0
$data[27,6] = 153
$data[27,7] = 4172
$data[28,0] = 'Chequia'
$data[28,1] = 223065
$data[28,2] = 0
$data[28,3] = 87225
$data[28,4] = 133995
$data[28,5] = 0
$data[28,6] = 0
$data[28,7] = 1845
$data[29,0] = 'Canada'
$data[29,1] = 209974
$data[29,2] = 826
$data[29,3] = 176538
$data[29,4] = 23565
$data[29,5] = 0
$data[29,6] = 9
$data[29,7] = 9871
$data[30,0] = 'Rumania'
$data[30,1] = 201032
$data[30,2] = 5028
$data[30,3] = 144429
$data[30,4] = 50358
$data[30,5] = 0
$data[30,6] = 82
$data[30,7] = 6245
$data[31,0] = 'Marruecos'
$data[31,1] = 186731
$data[31,2] = 0
$data[31,3] = 154481
$data[31,4] = 29118
$data[31,5] = 0
$data[31,6] = 0
$data[31,7] = 3132
$data[32,0] = 'Ecuador'
$data[32,1] = 156451
$data[32,2] = 0
$data[32,3] = 134187
$data[32,4] = 9764
$data[32,5] = 0
$data[32,6] = 0
$data[32,7] = 12500
$data[33,0] = 'Nepal'
$data[33,1] = 153008
$data[33,2] = 4499
$data[33,3] = 105488
$data[33,4] = 46691
$data[33,5] = 0
$data[33,6] = 17
$data[33,7] = 829
$data[34,0] = 'Bolivia'
$data[34,1] = 140445
$data[34,2] = 217
$data[34,3] = 106950
$data[34,4] = 24911
$data[34,5] = 0
$data[34,6] = 26
$data[34,7] = 8584
$data[35,0] = 'Catar'
$data[35,1] = 130711
$data[35,2] = 249
$data[35,3] = 127599
$data[35,4] = 2883
$data[35,5] = 0
$data[35,6] = 1
$data[35,7] = 229
$data[36,0] = 'Panama'
$data[36,1] = 127227
$data[36,2] = 0
$data[36,3] = 103398
$data[36,4] = 21217
$data[36,5] = 0
$data[36,6] = 0
$data[36,7] = 2612
$data[37,0] = 'Republica Dominicana'
$data[37,1] = 122873
$data[37,2] = 0
$data[37,3] = 100920
$data[37,4] = 19741
$data[37,5] = 0
$data[37,6] = 0
$data[37,7] = 2212
$data[38,0] = 'Emiratos Arabes Unidos'
$data[38,1] = 122273
$data[38,2] = 1563
$data[38,3] = 115068
$data[38,4] = 6730
$data[38,5] = 0
$data[38,6] = 1
$data[38,7] = 475
$data[39,0] = 'Kuwait'
$data[39,1] = 120232
$data[39,2] = 812
$data[39,3] = 111440
$data[39,4] = 8052
$data[39,5] = 0
$data[39,6] = 10
$data[39,7] = 740
$data[40,0] = 'Portugal'
$data[40,1] = 112440
$data[40,2] = 2899
$data[40,3] = 65880
$data[40,4] = 44284
$data[40,5] = 0
$data[40,6] = 31
$data[40,7] = 2276
$data[41,0] = 'Oman'
$data[41,1] = 111837
$data[41,2] = 0
$data[41,3] = 97949
$data[41,4] = 12741
$data[41,5] = 0
$data[41,6] = 0
$data[41,7] = 1147
$data[42,0] = 'Suecia'
$data[42,1] = 110594
$data[42,2] = 0
$data[42,3] = 0
$data[42,4] = 0
$data[42,5] = 0
$data[42,6] = 7
$data[42,7] = 5933
$data[43,0] = 'Kazajistan'
$data[43,1] = 110086
$data[43,2] = 179
$data[43,3] = 105493
$data[43,4] = 2797
$data[43,5] = 0
$data[43,6] = 0
$data[43,7] = 1796
$data[44,0] = 'Egipto'
$data[44,1] = 106060
$data[44,2] = 0
$data[44,3] = 98624
$data[44,4] = 1270
$data[44,5] = 0
$data[44,6] = 0
$data[44,7] = 6166
$data[45,0] = 'Suiza'
$data[45,1] = 103653
$data[45,2] = 6634
$data[45,3] = 55800
$data[45,4] = 45790
$data[45,5] = 0
$data[45,6] = 11
$data[45,7] = 2063
$data[46,0] = 'Guatemala'
$data[46,1] = 103172
$data[46,2] = 0
$data[46,3] = 92665
$data[46,4] = 6927
$data[46,5] = 0
$data[46,6] = 0
$data[46,7] = 3580
$data[47,0] = 'Costa Rica'
$data[47,1] = 100616
$data[47,2] = 0
$data[47,3] = 61162
$data[47,4] = 38203
$data[47,5] = 0
$data[47,6] = 0
$data[47,7] = 1251
$data[48,0] = 'Japon'
$data[48,1] = 95138
$data[48,2] = 614
$data[48,3] = 88245
$data[48,4] = 5199
$data[48,5] = 0
$data[48,6] = 9
$data[48,7] = 1694
$data[49,0] = 'Etiopia'
$data[49,1] = 91693
$data[49,2] = 0
$data[49,3] = 45260
$data[49,4] = 45037
$data[49,5] = 0
$data[49,6] = 0
$data[49,7] = 1396
$data[50,0] = 'Honduras'
$data[50,1] = 91509
$data[50,2] = 431
$data[50,3] = 37132
$data[50,4] = 51773
$data[50,5] = 0
$data[50,6] = 8
$data[50,7] = 2604
$data[51,0] = 'Bielorrusia'
$data[51,1] = 91167
$data[51,2] = 787
$data[51,3] = 82136
$data[51,4] = 8086
$data[51,5] = 0
$data[51,6] = 0
$data[51,7] = 945
$data[52,0] = 'Venezuela'
$data[52,1] = 88416
$data[52,2] = 0
$data[52,3] = 82284
$data[52,4] = 5373
$data[52,5] = 0
$data[52,6] = 0
$data[52,7] = 759
$data[53,0] = 'China'
$data[53,1] = 85747
$data[53,2] = 18
$data[53,3] = 80865
$data[53,4] = 248
$data[53,5] = 0
$data[53,6] = 0
$data[53,7] = 4634
$data[54,0] = 'Barein'
$data[54,1] = 79211
$data[54,2] = 0
$data[54,3] = 75840
$data[54,4] = 3061
$data[54,5] = 0
$data[54,6] = 2
$data[54,7] = 310
$data[55,0] = 'Austria'
$data[55,1] = 74415
$data[55,2] = 2571
$data[55,3] = 55195
$data[55,4] = 18266
$data[55,5] = 0
$data[55,6] = 13
$data[55,7] = 954
$data[56,0] = 'Armenia'
$data[56,1] = 73310
$data[56,2] = 2474
$data[56,3] = 50276
$data[56,4] = 21889
$data[56,5] = 0
$data[56,6] = 14
$data[56,7] = 1145
$data[57,0] = 'Moldavia'
$data[57,1] = 70256
$data[57,2] = 688
$data[57,3] = 51102
$data[57,4] = 17500
$data[57,5] = 0
$data[57,6] = 13
$data[57,7] = 1654
$data[58,0] = 'Libano'
$data[58,1] = 67027
$data[58,2] = 0
$data[58,3] = 31409
$data[58,4] = 35066
$data[58,5] = 0
$data[58,6] = 0
$data[58,7] = 552
$data[59,0] = 'Uzbekistan'
$data[59,1] = 64633
$data[59,2] = 194
$data[59,3] = 61734
$data[59,4] = 2358
$data[59,5] = 0
$data[59,6] = 1
$data[59,7] = 541
$data[60,0] = 'Nigeria'
$data[60,1] = 61805
$data[60,2] = 0
$data[60,3] = 56985
$data[60,4] = 3693
$data[60,5] = 0
$data[60,6] = 0
$data[60,7] = 1127
$data[61,0] = 'Singapur'
$data[61,1] = 57951
$data[61,2] = 10
$data[61,3] = 57829
$data[61,4] = 94
$data[61,5] = 0
$data[61,6] = 0
$data[61,7] = 28
$data[62,0] = 'Paraguay'
$data[62,1] = 57526
$data[62,2] = 0
$data[62,3] = 38187
$data[62,4] = 18077
$data[62,5] = 0
$data[62,6] = 0
$data[62,7] = 1262
$data[63,0] = 'Argelia'
$data[63,1] = 55357
$data[63,2] = 0
$data[63,3] = 38618
$data[63,4] = 14851
$data[63,5] = 0
$data[63,6] = 0
$data[63,7] = 1888
$data[64,0] = 'Kirguistan'
$data[64,1] = 54588
$data[64,2] = 582
$data[64,3] = 47050
$data[64,4] = 6412
$data[64,5] = 0
$data[64,6] = 4
$data[64,7] = 1126
$data[65,0] = 'Irlanda'
$data[65,1] = 54476
$data[65,2] = 0
$data[65,3] = 23364
$data[65,4] = 29241
$data[65,5] = 0
$data[65,6] = 0
$data[65,7] = 1871
$data[66,0] = 'Hungria'
$data[66,1] = 54278
$data[66,2] = 2066
$data[66,3] = 15655
$data[66,4] = 37271
$data[66,5] = 0
$data[66,6] = 47
$data[66,7] = 1352
$data[67,0] = 'Libia'
$data[67,1] = 53384
$data[67,2] = 764
$data[67,3] = 29619
$data[67,4] = 22991
$data[67,5] = 0
$data[67,6] = 6
$data[67,7] = 774
$data[68,0] = 'Estado de Palestina'
$data[68,1] = 49579
$data[68,2] = 445
$data[68,3] = 42850
$data[68,4] = 6290
$data[68,5] = 0
$data[68,6] = 4
$data[68,7] = 439
$data[69,0] = 'Azerbaiyan'
$data[69,1] = 48221
$data[69,2] = 803
$data[69,3] = 40831
$data[69,4] = 6734
$data[69,5] = 0
$data[69,6] = 8
$data[69,7] = 656
$data[70,0] = 'Ghana'
$data[70,1] = 47601
$data[70,2] = 63
$data[70,3] = 46824
$data[70,4] = 463
$data[70,5] = 0
$data[70,6] = 2
$data[70,7] = 314
$data[71,0] = 'Tunez'
$data[71,1] = 47214
$data[71,2] = 1322
$data[71,3] = 5032
$data[71,4] = 41398
$data[71,5] = 0
$data[71,6] = 44
$data[71,7] = 784
$data[72,0] = 'Kenia'
$data[72,1] = 47212
$data[72,2] = 0
$data[72,3] = 33050
$data[72,4] = 13292
$data[72,5] = 0
$data[72,6] = 0
$data[72,7] = 870
$data[73,0] = 'Jordania'
$data[73,1] = 46441
$data[73,2] = 0
$data[73,3] = 7340
$data[73,4] = 38620
$data[73,5] = 0
$data[73,6] = 0
$data[73,7] = 481
$data[74,0] = 'Birmania'
$data[74,1] = 42365
$data[74,2] = 1357
$data[74,3] = 22445
$data[74,4] = 18882
$data[74,5] = 0
$data[74,6] = 33
$data[74,7] = 1038
$data[75,0] = 'Afganistan'
$data[75,1] = 40687
$data[75,2] = 61
$data[75,3] = 34010
$data[75,4] = 5170
$data[75,5] = 0
$data[75,6] = 2
$data[75,7] = 1507
$data[76,0] = 'Dinamarca'
$data[76,1] = 38622
$data[76,2] = 859
$data[76,3] = 31295
$data[76,4] = 6630
$data[76,5] = 0
$data[76,6] = 3
$data[76,7] = 697
$data[77,0] = 'Bosnia y Herzegovina'
$data[77,1] = 38493
$data[77,2] = 1179
$data[77,3] = 26260
$data[77,4] = 11168
$data[77,5] = 0
$data[77,6] = 14
$data[77,7] = 1065
$data[78,0] = 'Serbia'
$data[78,1] = 38115
$data[78,2] = 579
$data[78,3] = 31536
$data[78,4] = 5793
$data[78,5] = 0
$data[78,6] = 3
$data[78,7] = 786
$data[79,0] = 'Eslovaquia'
$data[79,1] = 37911
$data[79,2] = 2581
$data[79,3] = 8859
$data[79,4] = 28918
$data[79,5] = 0
$data[79,6] = 19
$data[79,7] = 134
$data[80,0] = 'Bulgaria'
$data[80,1] = 34930
$data[80,2] = 0
$data[80,3] = 17833
$data[80,4] = 16033
$data[80,5] = 0
$data[80,6] = 0
$data[80,7] = 1064
$data[81,0] = 'El Salvador'
$data[81,1] = 32421
$data[81,2] = 159
$data[81,3] = 28127
$data[81,4] = 3354
$data[81,5] = 0
$data[81,6] = 4
$data[81,7] = 940
$data[82,0] = 'Croacia'
$data[82,1] = 31717
$data[82,2] = 1867
$data[82,3] = 22910
$data[82,4] = 8394
$data[82,5] = 0
$data[82,6] = 7
$data[82,7] = 413
$data[83,0] = 'Grecia'
$data[83,1] = 28216
$data[83,2] = 0
$data[83,3] = 9989
$data[83,4] = 17678
$data[83,5] = 0
$data[83,6] = 0
$data[83,7] = 549
$data[84,0] = 'Australia'
$data[84,1] = 27484
$data[84,2] = 18
$data[84,3] = 25169
$data[84,4] = 1410
$data[84,5] = 0
$data[84,6] = 0
$data[84,7] = 905
$data[85,0] = 'Corea del Sur'
$data[85,1] = 25698
$data[85,2] = 155
$data[85,3] = 23717
$data[85,4] = 1526
$data[85,5] = 0
$data[85,6] = 2
$data[85,7] = 455
$data[86,0] = 'Republica de Macedonia'
$data[86,1] = 25473
$data[86,2] = 0
$data[86,3] = 18047
$data[86,4] = 6552
$data[86,5] = 0
$data[86,6] = 0
$data[86,7] = 874
$data[87,0] = 'Georgia'
$data[87,1] = 24562
$data[87,2] = 1759
$data[87,3] = 9751
$data[87,4] = 14628
$data[87,5] = 0
$data[87,6] = 5
$data[87,7] = 183
$data[88,0] = 'Malasia'
$data[88,1] = 24514
$data[88,2] = 710
$data[88,3] = 15884
$data[88,4] = 8416
$data[88,5] = 0
$data[88,6] = 10
$data[88,7] = 214
$data[89,0] = 'Camerun'
$data[89,1] = 21570
$data[89,2] = 0
$data[89,3] = 20117
$data[89,4] = 1028
$data[89,5] = 0
$data[89,6] = 0
$data[89,7] = 425
$data[90,0] = 'Costa de Marfil'
$data[90,1] = 20390
$data[90,2] = 0
$data[90,3] = 20088
$data[90,4] = 181
$data[90,5] = 0
$data[90,6] = 0
$data[90,7] = 121
$data[91,0] = 'Eslovenia'
$data[91,1] = 19307
$data[91,2] = 1656
$data[91,3] = 7659
$data[91,4] = 11434
$data[91,5] = 0
$data[91,6] = 3
$data[91,7] = 214
$data[92,0] = 'Albania'
$data[92,1] = 18556
$data[92,2] = 306
$data[92,3] = 10466
$data[92,4] = 7621
$data[92,5] = 0
$data[92,6] = 4
$data[92,7] = 469
$data[93,0] = 'Noruega'
$data[93,1] = 17405
$data[93,2] = 173
$data[93,3] = 11863
$data[93,4] = 5263
$data[93,5] = 0
$data[93,6] = 0
$data[93,7] = 279
$data[94,0] = 'Madagascar'
$data[94,1] = 16810
$data[94,2] = 0
$data[94,3] = 16215
$data[94,4] = 357
$data[94,5] = 0
$data[94,6] = 0
$data[94,7] = 238
$data[95,0] = 'Montenegro'
$data[95,1] = 16259
$data[95,2] = 0
$data[95,3] = 12093
$data[95,4] = 3913
$data[95,5] = 0
$data[95,6] = 0
$data[95,7] = 253
$data[96,0] = 'Zambia'
$data[96,1] = 16095
$data[96,2] = 60
$data[96,3] = 15179
$data[96,4] = 570
$data[96,5] = 0
$data[96,6] = 0
$data[96,7] = 346
$data[97,0] = 'Senegal'
$data[97,1] = 15525
$data[97,2] = 17
$data[97,3] = 14082
$data[97,4] = 1122
$data[97,5] = 0
$data[97,6] = 0
$data[97,7] = 321
$data[98,0] = 'Finlandia'
$data[98,1] = 14474
$data[98,2] = 219
$data[98,3] = 9800
$data[98,4] = 4321
$data[98,5] = 0
$data[98,6] = 0
$data[98,7] = 353
$data[99,0] = 'Sudan'
$data[99,1] = 13724
$data[99,2] = 0
$data[99,3] = 6764
$data[99,4] = 6124
$data[99,5] = 0
$data[99,6] = 0
$data[99,7] = 836
$data[100,0] = 'Namibia'
$data[100,1] = 12460
$data[100,2] = 0
$data[100,3] = 10609
$data[100,4] = 1718
$data[100,5] = 0
$data[100,6] = 0
$data[100,7] = 133
$data[101,0] = 'Luxemburgo'
$data[101,1] = 12333
$data[101,2] = 0
$data[101,3] = 8474
$data[101,4] = 3719
$data[101,5] = 0
$data[101,6] = 0
$data[101,7] = 140
$data[102,0] = 'Guinea'
$data[102,1] = 11635
$data[102,2] = 0
$data[102,3] = 10474
$data[102,4] = 1090
$data[102,5] = 0
$data[102,6] = 0
$data[102,7] = 71
$data[103,0] = 'Mozambique'
$data[103,1] = 11559
$data[103,2] = 0
$data[103,3] = 9226
$data[103,4] = 2252
$data[103,5] = 0
$data[103,6] = 0
$data[103,7] = 81
$data[104,0] = 'Maldivas'
$data[104,1] = 11358
$data[104,2] = 0
$data[104,3] = 10383
$data[104,4] = 938
$data[104,5] = 0
$data[104,6] = 0
$data[104,7] = 37
$data[105,0] = 'Uganda'
$data[105,1] = 11163
$data[105,2] = 122
$data[105,3] = 7269
$data[105,4] = 3795
$data[105,5] = 0
$data[105,6] = 1
$data[105,7] = 99
$data[106,0] = 'Consejo Danes para los Refugiados'
$data[106,1] = 11122
$data[106,2] = 25
$data[106,3] = 10379
$data[106,4] = 439
$data[106,5] = 0
$data[106,6] = 0
$data[106,7] = 304
$data[107,0] = 'Tayikistan'
$data[107,1] = 10695
$data[107,2] = 42
$data[107,3] = 9782
$data[107,4] = 832
$data[107,5] = 0
$data[107,6] = 0
$data[107,7] = 81
$data[108,0] = 'Guayana Francesa'
$data[108,1] = 10342
$data[108,2] = 0
$data[108,3] = 9995
$data[108,4] = 278
$data[108,5] = 0
$data[108,6] = 0
$data[108,7] = 69
$data[109,0] = 'Lituania'
$data[109,1] = 9104
$data[109,2] = 442
$data[109,3] = 3978
$data[109,4] = 5000
$data[109,5] = 0
$data[109,6] = 1
$data[109,7] = 126
$data[110,0] = 'Haiti'
$data[110,1] = 9007
$data[110,2] = 0
$data[110,3] = 7311
$data[110,4] = 1465
$data[110,5] = 0
$data[110,6] = 0
$data[110,7] = 231
$data[111,0] = 'Gabon'
$data[111,1] = 8901
$data[111,2] = 0
$data[111,3] = 8479
$data[111,4] = 368
$data[111,5] = 0
$data[111,6] = 0
$data[111,7] = 54
$data[112,0] = 'Jamaica'
$data[112,1] = 8638
$data[112,2] = 38
$data[112,3] = 4156
$data[112,4] = 4300
$data[112,5] = 0
$data[112,6] = 3
$data[112,7] = 182
$data[113,0] = 'Angola'
$data[113,1] = 8582
$data[113,2] = 0
$data[113,3] = 3305
$data[113,4] = 5017
$data[113,5] = 0
$data[113,6] = 0
$data[113,7] = 260
$data[114,0] = 'Zimbabue'
$data[114,1] = 8242
$data[114,2] = 0
$data[114,3] = 7742
$data[114,4] = 264
$data[114,5] = 0
$data[114,6] = 0
$data[114,7] = 236
$data[115,0] = 'Cabo Verde'
$data[115,1] = 8122
$data[115,2] = 0
$data[115,3] = 6940
$data[115,4] = 1091
$data[115,5] = 0
$data[115,6] = 0
$data[115,7] = 91
$data[116,0] = 'Mauritania'
$data[116,1] = 7650
$data[116,2] = 0
$data[116,3] = 7369
$data[116,4] = 118
$data[116,5] = 0
$data[116,6] = 0
$data[116,7] = 163
$data[117,0] = 'Guadalupe'
$data[117,1] = 7329
$data[117,2] = 0
$data[117,3] = 2199
$data[117,4] = 5015
$data[117,5] = 0
$data[117,6] = 0
$data[117,7] = 115
$data[118,0] = 'Sri Lanka'
$data[118,1] = 6896
$data[118,2] = 609
$data[118,3] = 3644
$data[118,4] = 3238
$data[118,5] = 0
$data[118,6] = 0
$data[118,7] = 14
$data[119,0] = 'Cuba'
$data[119,1] = 6479
$data[119,2] = 58
$data[119,3] = 5899
$data[119,4] = 452
$data[119,5] = 0
$data[119,6] = 0
$data[119,7] = 128
$data[120,0] = 'Bahamas'
$data[120,1] = 6135
$data[120,2] = 0
$data[120,3] = 3705
$data[120,4] = 2300
$data[120,5] = 0
$data[120,6] = 0
$data[120,7] = 130
$data[121,0] = 'Botsuana'
$data[121,1] = 5923
$data[121,2] = 0
$data[121,3] = 927
$data[121,4] = 4975
$data[121,5] = 0
$data[121,6] = 0
$data[121,7] = 21
$data[122,0] = 'Malaui'
$data[122,1] = 5885
$data[122,2] = 11
$data[122,3] = 5287
$data[122,4] = 415
$data[122,5] = 0
$data[122,6] = 0
$data[122,7] = 183
$data[123,0] = 'Suazilandia'
$data[123,1] = 5814
$data[123,2] = 0
$data[123,3] = 5468
$data[123,4] = 230
$data[123,5] = 0
$data[123,6] = 0
$data[123,7] = 116
$data[124,0] = 'Republica de Yibuti'
$data[124,1] = 5522
$data[124,2] = 0
$data[124,3] = 5389
$data[124,4] = 72
$data[124,5] = 0
$data[124,6] = 0
$data[124,7] = 61
$data[125,0] = 'Trinidad yTobago'
$data[125,1] = 5446
$data[125,2] = 0
$data[125,3] = 3876
$data[125,4] = 1467
$data[125,5] = 0
$data[125,6] = 0
$data[125,7] = 103
$data[126,0] = 'Nicaragua'
$data[126,1] = 5434
$data[126,2] = 0
$data[126,3] = 4225
$data[126,4] = 1054
$data[126,5] = 0
$data[126,6] = 0
$data[126,7] = 155
$data[127,0] = 'Hong Kong'
$data[127,1] = 5285
$data[127,2] = 4
$data[127,3] = 5029
$data[127,4] = 151
$data[127,5] = 0
$data[127,6] = 0
$data[127,7] = 105
$data[128,0] = 'Siria'
$data[128,1] = 5267
$data[128,2] = 0
$data[128,3] = 1655
$data[128,4] = 3352
$data[128,5] = 0
$data[128,6] = 0
$data[128,7] = 260
$data[129,0] = 'Malta'
$data[129,1] = 5258
$data[129,2] = 121
$data[129,3] = 3439
$data[129,4] = 1770
$data[129,5] = 0
$data[129,6] = 0
$data[129,7] = 49
$data[130,0] = 'Polinesia Francesa'
$data[130,1] = 5161
$data[130,2] = 0
$data[130,3] = 3536
$data[130,4] = 1606
$data[130,5] = 0
$data[130,6] = 0
$data[130,7] = 19
$data[131,0] = 'Congo'
$data[131,1] = 5156
$data[131,2] = 0
$data[131,3] = 3887
$data[131,4] = 1177
$data[131,5] = 0
$data[131,6] = 0
$data[131,7] = 92
$data[132,0] = 'Surinam'
$data[132,1] = 5154
$data[132,2] = 0
$data[132,3] = 4995
$data[132,4] = 50
$data[132,5] = 0
$data[132,6] = 0
$data[132,7] = 109
$data[133,0] = 'Reunion'
$data[133,1] = 5149
$data[133,2] = 134
$data[133,3] = 4630
$data[133,4] = 499
$data[133,5] = 0
$data[133,6] = 1
$data[133,7] = 20
$data[134,0] = 'Guinea Ecuatorial'
$data[134,1] = 5074
$data[134,2] = 0
$data[134,3] = 4961
$data[134,4] = 30
$data[134,5] = 0
$data[134,6] = 0
$data[134,7] = 83
$data[135,0] = 'Ruanda'
$data[135,1] = 5017
$data[135,2] = 0
$data[135,3] = 4803
$data[135,4] = 180
$data[135,5] = 0
$data[135,6] = 0
$data[135,7] = 34
$data[136,0] = 'Republica de Africa Central'
$data[136,1] = 4862
$data[136,2] = 0
$data[136,3] = 1924
$data[136,4] = 2876
$data[136,5] = 0
$data[136,6] = 0
$data[136,7] = 62
$data[137,0] = 'Aruba'
$data[137,1] = 4389
$data[137,2] = 0
$data[137,3] = 4120
$data[137,4] = 233
$data[137,5] = 0
$data[137,6] = 0
$data[137,7] = 36
$data[138,0] = 'Islandia'
$data[138,1] = 4308
$data[138,2] = 40
$data[138,3] = 3187
$data[138,4] = 1110
$data[138,5] = 0
$data[138,6] = 0
$data[138,7] = 11
$data[139,0] = 'Estonia'
$data[139,1] = 4300
$data[139,2] = 53
$data[139,3] = 3418
$data[139,4] = 809
$data[139,5] = 0
$data[139,6] = 2
$data[139,7] = 73
$data[140,0] = 'Letonia'
$data[140,1] = 4208
$data[140,2] = 250
$data[140,3] = 1357
$data[140,4] = 2801
$data[140,5] = 0
$data[140,6] = 1
$data[140,7] = 50
$data[141,0] = 'Mayotte'
$data[141,1] = 4203
$data[141,2] = 0
$data[141,3] = 2964
$data[141,4] = 1195
$data[141,5] = 0
$data[141,6] = 0
$data[141,7] = 44
$data[142,0] = 'Somalia'
$data[142,1] = 3897
$data[142,2] = 0
$data[142,3] = 3166
$data[142,4] = 629
$data[142,5] = 0
$data[142,6] = 0
$data[142,7] = 102
$data[143,0] = 'Guyana'
$data[143,1] = 3877
$data[143,2] = 0
$data[143,3] = 2853
$data[143,4] = 907
$data[143,5] = 0
$data[143,6] = 0
$data[143,7] = 117
$data[144,0] = 'Principado de Andorra'
$data[144,1] = 3811
$data[144,2] = 0
$data[144,3] = 2470
$data[144,4] = 1278
$data[144,5] = 0
$data[144,6] = 0
$data[144,7] = 63
$data[145,0] = 'Tailandia'
$data[145,1] = 3727
$data[145,2] = 8
$data[145,3] = 3518
$data[145,4] = 150
$data[145,5] = 0
$data[145,6] = 0
$data[145,7] = 59
$data[146,0] = 'Gambia'
$data[146,1] = 3659
$data[146,2] = 0
$data[146,3] = 2660
$data[146,4] = 880
$data[146,5] = 0
$data[146,6] = 0
$data[146,7] = 119
$data[147,0] = 'Mali'
$data[147,1] = 3440
$data[147,2] = 0
$data[147,3] = 2608
$data[147,4] = 700
$data[147,5] = 0
$data[147,6] = 0
$data[147,7] = 132
$data[148,0] = 'Republica de Chipre'
$data[148,1] = 3154
$data[148,2] = 0
$data[148,3] = 1444
$data[148,4] = 1685
$data[148,5] = 0
$data[148,6] = 0
$data[148,7] = 25
$data[149,0] = 'Belice'
$data[149,1] = 2995
$data[149,2] = 58
$data[149,3] = 1826
$data[149,4] = 1123
$data[149,5] = 0
$data[149,6] = 0
$data[149,7] = 46
$data[150,0] = 'Sudan del Sur'
$data[150,1] = 2872
$data[150,2] = 0
$data[150,3] = 1290
$data[150,4] = 1527
$data[150,5] = 0
$data[150,6] = 0
$data[150,7] = 55
$data[151,0] = 'Uruguay'
$data[151,1] = 2701
$data[151,2] = 0
$data[151,3] = 2204
$data[151,4] = 444
$data[151,5] = 0
$data[151,6] = 0
$data[151,7] = 53
$data[152,0] = 'Benin'
$data[152,1] = 2557
$data[152,2] = 0
$data[152,3] = 2330
$data[152,4] = 186
$data[152,5] = 0
$data[152,6] = 0
$data[152,7] = 41
$data[153,0] = 'Burkina Faso'
$data[153,1] = 2414
$data[153,2] = 0
$data[153,3] = 1869
$data[153,4] = 480
$data[153,5] = 0
$data[153,6] = 0
$data[153,7] = 65
$data[154,0] = 'Guinea-Bisau'
$data[154,1] = 2403
$data[154,2] = 0
$data[154,3] = 1818
$data[154,4] = 544
$data[154,5] = 0
$data[154,6] = 0
$data[154,7] = 41
$data[155,0] = 'Sierra Leona'
$data[155,1] = 2340
$data[155,2] = 0
$data[155,3] = 1777
$data[155,4] = 490
$data[155,5] = 0
$data[155,6] = 0
$data[155,7] = 73
$data[156,0] = 'Martinica'
$data[156,1] = 2257
$data[156,2] = 0
$data[156,3] = 98
$data[156,4] = 2135
$data[156,5] = 0
$data[156,6] = 0
$data[156,7] = 24
$data[157,0] = 'Togo'
$data[157,1] = 2139
$data[157,2] = 0
$data[157,3] = 1574
$data[157,4] = 513
$data[157,5] = 0
$data[157,6] = 0
$data[157,7] = 52
$data[158,0] = 'Yemen'
$data[158,1] = 2057
$data[158,2] = 0
$data[158,3] = 1344
$data[158,4] = 116
$data[158,5] = 0
$data[158,6] = 0
$data[158,7] = 597
$data[159,0] = 'Lesoto'
$data[159,1] = 1934
$data[159,2] = 11
$data[159,3] = 961
$data[159,4] = 930
$data[159,5] = 0
$data[159,6] = 0
$data[159,7] = 43
$data[160,0] = 'Nueva Zelanda'
$data[160,1] = 1923
$data[160,2] = 9
$data[160,3] = 1832
$data[160,4] = 66
$data[160,5] = 0
$data[160,6] = 0
$data[160,7] = 25
$data[161,0] = 'Republica del Chad'
$data[161,1] = 1410
$data[161,2] = 0
$data[161,3] = 1223
$data[161,4] = 91
$data[161,5] = 0
$data[161,6] = 0
$data[161,7] = 96
$data[162,0] = 'Liberia'
$data[162,1] = 1385
$data[162,2] = 0
$data[162,3] = 1278
$data[162,4] = 25
$data[162,5] = 0
$data[162,6] = 0
$data[162,7] = 82
$data[163,0] = 'Niger'
$data[163,1] = 1215
$data[163,2] = 0
$data[163,3] = 1128
$data[163,4] = 18
$data[163,5] = 0
$data[163,6] = 0
$data[163,7] = 69
$data[164,0] = 'Vietnam'
$data[164,1] = 1148
$data[164,2] = 0
$data[164,3] = 1049
$data[164,4] = 64
$data[164,5] = 0
$data[164,6] = 0
$data[164,7] = 35
$data[165,0] = 'Santo Tome y Principe'
$data[165,1] = 935
$data[165,2] = 0
$data[165,3] = 898
$data[165,4] = 22
$data[165,5] = 0
$data[165,6] = 0
$data[165,7] = 15
$data[166,0] = 'San Marino'
$data[166,1] = 819
$data[166,2] = 17
$data[166,3] = 716
$data[166,4] = 61
$data[166,5] = 0
$data[166,6] = 0
$data[166,7] = 42
$data[167,0] = 'Curazao'
$data[167,1] = 804
$data[167,2] = 0
$data[167,3] = 509
$data[167,4] = 294
$data[167,5] = 0
$data[167,6] = 0
$data[167,7] = 1
$data[168,0] = 'San Martin (Parte Holandesa)'
$data[168,1] = 769
$data[168,2] = 0
$data[168,3] = 681
$data[168,4] = 66
$data[168,5] = 0
$data[168,6] = 0
$data[168,7] = 22
$data[169,0] = 'Crucero'
$data[169,1] = 712
$data[169,2] = 0
$data[169,3] = 659
$data[169,4] = 40
$data[169,5] = 0
$data[169,6] = 0
$data[169,7] = 13
$data[170,0] = 'Islas Turcas y Caicos'
$data[170,1] = 698
$data[170,2] = 0
$data[170,3] = 689
$data[170,4] = 3
$data[170,5] = 0
$data[170,6] = 0
$data[170,7] = 6
$data[171,0] = 'Gibraltar'
$data[171,1] = 641
$data[171,2] = 11
$data[171,3] = 500
$data[171,4] = 141
$data[171,5] = 0
$data[171,6] = 0
$data[171,7] = 0
$data[172,0] = 'Papua Nueva Guinea'
$data[172,1] = 583
$data[172,2] = 0
$data[172,3] = 545
$data[172,4] = 31
$data[172,5] = 0
$data[172,6] = 0
$data[172,7] = 7
$data[173,0] = 'Burundi'
$data[173,1] = 551
$data[173,2] = 0
$data[173,3] = 497
$data[173,4] = 53
$data[173,5] = 0
$data[173,6] = 0
$data[173,7] = 1
$data[174,0] = 'Taiwan'
$data[174,1] = 548
$data[174,2] = 0
$data[174,3] = 497
$data[174,4] = 44
$data[174,5] = 0
$data[174,6] = 0
$data[174,7] = 7
$data[175,0] = 'San Martin (Parte Francesa)'
$data[175,1] = 538
$data[175,2] = 0
$data[175,3] = 422
$data[175,4] = 108
$data[175,5] = 0
$data[175,6] = 0
$data[175,7] = 8
$data[176,0] = 'Comoras'
$data[176,1] = 517
$data[176,2] = 0
$data[176,3] = 494
$data[176,4] = 16
$data[176,5] = 0
$data[176,6] = 0
$data[176,7] = 7
$data[177,0] = 'Tanzania'
$data[177,1] = 509
$data[177,2] = 0
$data[177,3] = 183
$data[177,4] = 305
$data[177,5] = 0
$data[177,6] = 0
$data[177,7] = 21
$data[178,0] = 'Islas Feroe'
$data[178,1] = 490
$data[178,2] = 0
$data[178,3] = 473
$data[178,4] = 17
$data[178,5] = 0
$data[178,6] = 0
$data[178,7] = 0
$data[179,0] = 'Eritrea'
$data[179,1] = 457
$data[179,2] = 0
$data[179,3] = 391
$data[179,4] = 66
$data[179,5] = 0
$data[179,6] = 0
$data[179,7] = 0
$data[180,0] = 'Mauricio'
$data[180,1] = 425
$data[180,2] = 0
$data[180,3] = 386
$data[180,4] = 29
$data[180,5] = 0
$data[180,6] = 0
$data[180,7] = 10
$data[181,0] = 'Isla de Man'
$data[181,1] = 348
$data[181,2] = 0
$data[181,3] = 321
$data[181,4] = 3
$data[181,5] = 0
$data[181,6] = 0
$data[181,7] = 24
$data[182,0] = 'Butan'
$data[182,1] = 336
$data[182,2] = 4
$data[182,3] = 306
$data[182,4] = 30
$data[182,5] = 0
$data[182,6] = 0
$data[182,7] = 0
$data[183,0] = 'Mongolia'
$data[183,1] = 328
$data[183,2] = 0
$data[183,3] = 312
$data[183,4] = 16
$data[183,5] = 0
$data[183,6] = 0
$data[183,7] = 0
$data[184,0] = 'Liechtenstein'
$data[184,1] = 324
$data[184,2] = 42
$data[184,3] = 170
$data[184,4] = 153
$data[184,5] = 0
$data[184,6] = 0
$data[184,7] = 1
$data[185,0] = 'Camboya'
$data[185,1] = 286
$data[185,2] = 0
$data[185,3] = 280
$data[185,4] = 6
$data[185,5] = 0
$data[185,6] = 0
$data[185,7] = 0
$data[186,0] = 'Monaco'
$data[186,1] = 281
$data[186,2] = 0
$data[186,3] = 233
$data[186,4] = 46
$data[186,5] = 0
$data[186,6] = 0
$data[186,7] = 2
$data[187,0] = 'Islas Caimanes'
$data[187,1] = 236
$data[187,2] = 0
$data[187,3] = 215
$data[187,4] = 20
$data[187,5] = 0
$data[187,6] = 0
$data[187,7] = 1
$data[188,0] = 'Barbados'
$data[188,1] = 224
$data[188,2] = 0
$data[188,3] = 207
$data[188,4] = 10
$data[188,5] = 0
$data[188,6] = 0
$data[188,7] = 7
$data[189,0] = 'Bermudas'
$data[189,1] = 188
$data[189,2] = 0
$data[189,3] = 175
$data[189,4] = 4
$data[189,5] = 0
$data[189,6] = 0
$data[189,7] = 9
$data[190,0] = 'Seychelles'
$data[190,1] = 153
$data[190,2] = 2
$data[190,3] = 149
$data[190,4] = 4
$data[190,5] = 0
$data[190,6] = 0
$data[190,7] = 0
$data[191,0] = 'Bonaire, San Eustaquio y Saba'
$data[191,1] = 150
$data[191,2] = 0
$data[191,3] = 121
$data[191,4] = 26
$data[191,5] = 0
$data[191,6] = 0
$data[191,7] = 3
$data[192,0] = 'Brunei'
$data[192,1] = 148
$data[192,2] = 0
$data[192,3] = 143
$data[192,4] = 2
$data[192,5] = 0
$data[192,6] = 0
$data[192,7] = 3
$data[193,0] = 'Antigua y Barbuda'
$data[193,1] = 122
$data[193,2] = 0
$data[193,3] = 107
$data[193,4] = 12
$data[193,5] = 0
$data[193,6] = 0
$data[193,7] = 3
$data[194,0] = 'San Bartolome'
$data[194,1] = 77
$data[194,2] = 0
$data[194,3] = 66
$data[194,4] = 11
$data[194,5] = 0
$data[194,6] = 0
$data[194,7] = 0
$data[195,0] = 'Islas Virgenes Britanicas'
$data[195,1] = 71
$data[195,2] = 0
$data[195,3] = 70
$data[195,4] = 0
$data[195,5] = 0
$data[195,6] = 0
$data[195,7] = 1
$data[196,0] = 'San Vicente y las Granadinas'
$data[196,1] = 68
$data[196,2] = 0
$data[196,3] = 64
$data[196,4] = 4
$data[196,5] = 0
$data[196,6] = 0
$data[196,7] = 0
$data[197,0] = 'Macao'
$data[197,1] = 46
$data[197,2] = 0
$data[197,3] = 46
$data[197,4] = 0
$data[197,5] = 0
$data[197,6] = 0
$data[197,7] = 0
$data[198,0] = 'Santa Lucia'
$data[198,1] = 42
$data[198,2] = 0
$data[198,3] = 27
$data[198,4] = 15
$data[198,5] = 0
$data[198,6] = 0
$data[198,7] = 0
$data[199,0] = 'Puerto Rico'
$data[199,1] = 39
$data[199,2] = 0
$data[199,3] = 1
$data[199,4] = 36
$data[199,5] = 0
$data[199,6] = 0
$data[199,7] = 2
$data[200,0] = 'Dominica'
$data[200,1] = 33
$data[200,2] = 0
$data[200,3] = 29
$data[200,4] = 4
$data[200,5] = 0
$data[200,6] = 0
$data[200,7] = 0
$data[201,0] = 'Fiyi'
$data[201,1] = 33
$data[201,2] = 0
$data[201,3] = 30
$data[201,4] = 1
$data[201,5] = 0
$data[201,6] = 0
$data[201,7] = 2
$data[202,0] = 'Guam'
$data[202,1] = 32
$data[202,2] = 0
$data[202,3] = 0
$data[202,4] = 31
$data[202,5] = 0
$data[202,6] = 0
$data[202,7] = 1
$data[203,0] = 'Timor Oriental'
$data[203,1] = 29
$data[203,2] = 0
$data[203,3] = 28
$data[203,4] = 1
$data[203,5] = 0
$data[203,6] = 0
$data[203,7] = 0
$data[204,0] = 'Santa Sede'
$data[204,1] = 27
$data[204,2] = 0
$data[204,3] = 15
$data[204,4] = 12
$data[204,5] = 0
$data[204,6] = 0
$data[204,7] = 0
$data[205,0] = 'Granada'
$data[205,1] = 27
$data[205,2] = 0
$data[205,3] = 24
$data[205,4] = 3
$data[205,5] = 0
$data[205,6] = 0
$data[205,7] = 0
$data[206,0] = 'Nueva Caledonia'
$data[206,1] = 27
$data[206,2] = 0
$data[206,3] = 27
$data[206,4] = 0
$data[206,5] = 0
$data[206,6] = 0
$data[206,7] = 0
$data[207,0] = 'Laos'
$data[207,1] = 24
$data[207,2] = 0
$data[207,3] = 22
$data[207,4] = 2
$data[207,5] = 0
$data[207,6] = 0
$data[207,7] = 0
$data[208,0] = 'San Cristobal y Nieves'
$data[208,1] = 19
$data[208,2] = 0
$data[208,3] = 19
$data[208,4] = 0
$data[208,5] = 0
$data[208,6] = 0
$data[208,7] = 0
$data[209,0] = 'Islas Virgenes de los Estados Unidos'
$data[209,1] = 17
$data[209,2] = 0
$data[209,3] = 0
$data[209,4] = 17
$data[209,5] = 0
$data[209,6] = 0
$data[209,7] = 0
$data[210,0] = 'Groenlandia'
$data[210,1] = 17
$data[210,2] = 0
$data[210,3] = 16
$data[210,4] = 1
$data[210,5] = 0
$data[210,6] = 0
$data[210,7] = 0
$data[211,0] = 'San Pedro y Miquelon'
$data[211,1] = 16
$data[211,2] = 0
$data[211,3] = 12
$data[211,4] = 4
$data[211,5] = 0
$data[211,6] = 0
$data[211,7] = 0
$data[212,0] = 'Montserrat'
$data[212,1] = 13
$data[212,2] = 0
$data[212,3] = 12
$data[212,4] = 0
$data[212,5] = 0
$data[212,6] = 0
$data[212,7] = 1
$data[213,0] = 'Islas Malvinas'
$data[213,1] = 13
$data[213,2] = 0
$data[213,3] = 13
$data[213,4] = 0
$data[213,5] = 0
$data[213,6] = 0
$data[213,7] = 0
$data[214,0] = 'Sahara Occidental'
$data[214,1] = 10
$data[214,2] = 0
$data[214,3] = 8
$data[214,4] = 1
$data[214,5] = 0
$data[214,6] = 0
$data[214,7] = 1
$data[215,0] = 'Islas Salomon'
$data[215,1] = 4
$data[215,2] = 0
$data[215,3] = 3
$data[215,4] = 1
$data[215,5] = 0
$data[215,6] = 0
$data[215,7] = 0
$data[216,0] = 'Anguila'
$data[216,1] = 3
$data[216,2] = 0
$data[216,3] = 3
$data[216,4] = 0
$data[216,5] = 0
$data[216,6] = 0
$data[216,7] = 0
$data[217,0] = 'Wallis y Futuna'
$data[217,1] = 1
$data[217,2] = 0
$data[217,3] = 1
$data[217,4] = 0
$data[217,5] = 0
$data[217,6] = 0
$data[217,7] = 0

$ws.Range("A4:H221").Value = $data

